$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.667449
$ws.Range("N2").Value = 2.002347
$ws.Range("O2").Value = 0.05507918922684819
$ws.Range("P2").Value = 0.05743692635237437
$ws.Range("Q2").Value = 0.17734564896
$ws.Range("R2").Value = 1.59611084064
$ws.Range("S2").Value = 0.05507918922684819
$ws.Range("T2").Value = 0.05743692635237437

# Row 3
$ws.Range("O3").Value = 0.7558000582962523
$ws.Range("P3").Value = 0.7881530736897923
$ws.Range("S3").Value = 0.7558000582962523
$ws.Range("T3").Value = 0.7881530736897923

# Row 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4726029999999999
$ws.Range("N4").Value = 1.417809
$ws.Range("O4").Value = 0.0390001184602511
$ws.Range("P4").Value = 0.04066956981718631
$ws.Range("Q4").Value = 0.1255737677866667
$ws.Range("R4").Value = 1.13016391008
$ws.Range("S4").Value = 0.0390001184602511
$ws.Range("T4").Value = 0.04066956981718631

# Row 5
$ws.Range("M5").Value = 1.4922995
$ws.Range("N5").Value = 2.984599
$ws.Range("O5").Value = 0.1231474562755071
$ws.Range("P5").Value = 0.08561263005581461
$ws.Range("Q5").Value = 0.3965139258133333
$ws.Range("R5").Value = 2.37908355488
$ws.Range("S5").Value = 0.1231474562755071
$ws.Range("T5").Value = 0.08561263005581461

# Row 6
$ws.Range("M6").Value = 0.3268606666666667
$ws.Range("N6").Value = 0.9805820000000001
$ws.Range("O6").Value = 0.0269731777411414
$ws.Range("P6").Value = 0.02812780008483244
$ws.Range("Q6").Value = 0.08684905820444444
$ws.Range("R6").Value = 0.78164152384
$ws.Range("S6").Value = 0.0269731777411414
$ws.Range("T6").Value = 0.02812780008483244
